$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 373.4
$ws.Range("J17").Value = 307.7857
$ws.Range("L17").Value = 923.3571000000001
$ws.Range("N17").Value = -1259.3571
$ws.Range("H40").Value = 8123.5
$ws.Range("I40").Value = 1995
$ws.Range("J40").Value = 8999
$ws.Range("K40").Value = 1995
$ws.Range("L40").Value = 8999
$ws.Range("M40").Value = -1820
$ws.Range("N40").Value = -9349
$ws.Range("H70").Value = 5168
$ws.Range("J70").Value = 6162.5
$ws.Range("L70").Value = 18487.5
$ws.Range("N70").Value = -19027.5
$ws.Range("H73").Value = 5168
$ws.Range("J73").Value = 6162.5
$ws.Range("L73").Value = 18487.5
$ws.Range("N73").Value = -20359.5
$ws.Range("H87").Value = 19833.334
$ws.Range("J87").Value = 19833.334
$ws.Range("L87").Value = 19833.334
$ws.Range("N87").Value = -22329.334
$ws.Range("H90").Value = 19833.334
$ws.Range("J90").Value = 19833.334
$ws.Range("L90").Value = 59500.00199999999
$ws.Range("N90").Value = -71980.00199999999
$ws.Range("H98").Value = 2599.3635
$ws.Range("I98").Value = 1871.5555
$ws.Range("J98").Value = 5874.5
$ws.Range("K98").Value = 1871.5555
$ws.Range("L98").Value = 5874.5
$ws.Range("M98").Value = -373.5554999999999
$ws.Range("N98").Value = -8870.5
$ws.Range("H107").Value = 1500
$ws.Range("I107").Value = 1500
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 420
$ws.Range("N107").ClearContents()
$ws.Range("H111").Value = 2950.1333
$ws.Range("I111").Value = 1979.6666
$ws.Range("J111").Value = 4405.8335
$ws.Range("K111").Value = 5938.9998
$ws.Range("L111").Value = 13217.5005
$ws.Range("M111").Value = -2871.9998
$ws.Range("N111").Value = -19351.5005
$ws.Range("H122").Value = 2599.3635
$ws.Range("I122").Value = 1871.5555
$ws.Range("J122").Value = 5874.5
$ws.Range("K122").Value = 5614.666499999999
$ws.Range("L122").Value = 17623.5
$ws.Range("M122").Value = -3164.666499999999
$ws.Range("N122").Value = -22523.5
$ws.Range("H132").Value = 39101.215
$ws.Range("I132").Value = 2401.5454
$ws.Range("J132").Value = 173666.67
$ws.Range("K132").Value = 7204.6362
$ws.Range("L132").Value = 521000.01
$ws.Range("M132").Value = -4674.6362
$ws.Range("N132").Value = -526060.01

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 33847.453
$ws.Range("I2").Value = 56835.168
$ws.Range("J2").Value = 6262.2
$ws.Range("K2").Value = 56835.168
$ws.Range("L2").Value = 6262.2
$ws.Range("M2").Value = -56722.168
$ws.Range("N2").Value = -6488.2
$ws.Range("H45").Value = 5334.7144
$ws.Range("I45").Value = 4473.6
$ws.Range("K45").Value = 4473.6
$ws.Range("M45").Value = -4096.6
$ws.Range("H61").Value = 13607.77
$ws.Range("I61").Value = 4001.375
$ws.Range("K61").Value = 4001.375
$ws.Range("M61").Value = -3789.375
$ws.Range("H116").Value = 33847.453
$ws.Range("I116").Value = 56835.168
$ws.Range("J116").Value = 6262.2
$ws.Range("K116").Value = 56835.168
$ws.Range("L116").Value = 6262.2
$ws.Range("M116").Value = -54541.168
$ws.Range("N116").Value = -10850.2
$ws.Range("H122").Value = 5184.074
$ws.Range("I122").Value = 5069.4287
$ws.Range("K122").Value = 15208.2861
$ws.Range("M122").Value = -12758.2861
$ws.Range("H128").Value = 54714.5
$ws.Range("I128").Value = 45000
$ws.Range("K128").Value = 45000
$ws.Range("M128").Value = -40020
$ws.Range("H132").Value = 2002.5
$ws.Range("I132").Value = 2002.5
$ws.Range("K132").Value = 6007.5
$ws.Range("M132").Value = -3477.5
$ws.Range("H136").Value = 13607.77
$ws.Range("I136").Value = 4001.375
$ws.Range("K136").Value = 12004.125
$ws.Range("M136").Value = -9454.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 33847.453
$ws.Range("I3").Value = 56835.168
$ws.Range("J3").Value = 6262.2
$ws.Range("K3").Value = 56835.168
$ws.Range("L3").Value = 6262.2
$ws.Range("M3").Value = -56721.168
$ws.Range("N3").Value = -6490.2
$ws.Range("H20").Value = 11750.5
$ws.Range("I20").Value = 6801.4
$ws.Range("K20").Value = 6801.4
$ws.Range("M20").Value = -6554.4
$ws.Range("H86").Value = 20041668
$ws.Range("J86").Value = 1989.3
$ws.Range("L86").Value = 1989.3
$ws.Range("N86").Value = -4235.3
$ws.Range("H89").Value = 20041668
$ws.Range("J89").Value = 1989.3
$ws.Range("L89").Value = 9946.5
$ws.Range("N89").Value = -21178.5
$ws.Range("H105").Value = 1591.2222
$ws.Range("J105").Value = 3698.4285
$ws.Range("L105").Value = 3698.4285
$ws.Range("N105").Value = -7192.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1052
$ws.Range("I2").Value = 104
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 104
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = 9
$ws.Range("N2").Value = -2226
$ws.Range("H5").Value = 882.125
$ws.Range("I5").Value = 663.6667
$ws.Range("J5").Value = 1537.5
$ws.Range("K5").Value = 663.6667
$ws.Range("L5").Value = 1537.5
$ws.Range("M5").Value = -551.6667
$ws.Range("N5").Value = -1761.5
$ws.Range("H7").Value = 226.46153
$ws.Range("I7").Value = 84.90000000000001
$ws.Range("K7").Value = 84.90000000000001
$ws.Range("M7").Value = 28.09999999999999
$ws.Range("H10").Value = 1004
$ws.Range("J10").Value = 1308
$ws.Range("L10").Value = 1308
$ws.Range("N10").Value = -1586
$ws.Range("H11").Value = 3866.3333
$ws.Range("J11").Value = 3866.3333
$ws.Range("L11").Value = 3866.3333
$ws.Range("N11").Value = -4146.3333
$ws.Range("H12").Value = 995
$ws.Range("I12").Value = 995
$ws.Range("K12").Value = 995
$ws.Range("M12").Value = -825
$ws.Range("H44").Value = 163
$ws.Range("I44").Value = 163
$ws.Range("K44").Value = 163
$ws.Range("M44").Value = 279
$ws.Range("H58").Value = 2352.9614
$ws.Range("I58").Value = 1528.8334
$ws.Range("J58").Value = 3059.3572
$ws.Range("K58").Value = 1528.8334
$ws.Range("L58").Value = 3059.3572
$ws.Range("M58").Value = -1325.8334
$ws.Range("N58").Value = -3465.3572
$ws.Range("H132").Value = 3437.625
$ws.Range("I132").Value = 2897.4
$ws.Range("K132").Value = 8692.200000000001
$ws.Range("M132").Value = -6162.200000000001
$ws.Range("H134").Value = 8872
$ws.Range("I134").Value = 8872
$ws.Range("K134").Value = 26616
$ws.Range("M134").Value = -24081
$ws.Range("H136").Value = 2352.9614
$ws.Range("I136").Value = 1528.8334
$ws.Range("J136").Value = 3059.3572
$ws.Range("K136").Value = 4586.5002
$ws.Range("L136").Value = 9178.071599999999
$ws.Range("M136").Value = -2036.5002
$ws.Range("N136").Value = -14278.0716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 832.5454999999999
$ws.Range("I8").Value = 832.5454999999999
$ws.Range("K8").Value = 2497.6365
$ws.Range("M8").Value = -2358.6365
$ws.Range("H118").Value = 4499.75
$ws.Range("I118").Value = 5499
$ws.Range("J118").Value = 3500.5
$ws.Range("K118").Value = 16497
$ws.Range("L118").Value = 10501.5
$ws.Range("M118").Value = -15254
$ws.Range("N118").Value = -12987.5
$ws.Range("H138").Value = 5877.2
$ws.Range("I138").Value = 2921.7273
$ws.Range("J138").Value = 8199.357
$ws.Range("K138").Value = 8765.1819
$ws.Range("L138").Value = 24598.071
$ws.Range("M138").Value = -3625.1819
$ws.Range("N138").Value = -34878.071

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9192.821
$ws.Range("I132").Value = 8349.951999999999
$ws.Range("K132").Value = 25049.856
$ws.Range("M132").Value = -22519.856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 21
$ws.Range("I34").Value = 21
$ws.Range("K34").Value = 21
$ws.Range("M34").Value = 151
$ws.Range("H40").Value = 34950.375
$ws.Range("I40").Value = 11371.857
$ws.Range("K40").Value = 11371.857
$ws.Range("M40").Value = -11235.857
$ws.Range("H46").Value = 2777.3572
$ws.Range("H122").Value = 6190.2
$ws.Range("I122").Value = 6266.8887
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 18800.6661
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -16350.6661
$ws.Range("N122").Value = -21400
$ws.Range("H132").Value = 3083.923
$ws.Range("I132").Value = 2509.7
$ws.Range("K132").Value = 7529.099999999999
$ws.Range("M132").Value = -4999.099999999999
$ws.Range("H136").Value = 5020
$ws.Range("J136").Value = 5100
$ws.Range("L136").Value = 15300
$ws.Range("N136").Value = -20400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3849.5
$ws.Range("I62").Value = 3849
$ws.Range("J62").Value = 3850
$ws.Range("K62").Value = 3849
$ws.Range("L62").Value = 3850
$ws.Range("M62").Value = -3225
$ws.Range("N62").Value = -5098
$ws.Range("H65").Value = 3849.5
$ws.Range("I65").Value = 3849
$ws.Range("J65").Value = 3850
$ws.Range("K65").Value = 19245
$ws.Range("L65").Value = 19250
$ws.Range("M65").Value = -16125
$ws.Range("N65").Value = -25490
$ws.Range("H107").Value = 860.1429000000001
$ws.Range("I107").Value = 870.1667
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 2610.5001
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = -690.5001000000002
$ws.Range("N107").Value = -6240
$ws.Range("H132").Value = 5254.892
$ws.Range("I132").Value = 3860.9614
$ws.Range("J132").Value = 8549.637000000001
$ws.Range("K132").Value = 11582.8842
$ws.Range("L132").Value = 25648.911
$ws.Range("M132").Value = -9052.8842
$ws.Range("N132").Value = -30708.911
